$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (R5 / 1K resistor): quantity changes from 2 to 1, because R9 (row 7)
# is now broken out with its own quantity instead of being lumped into R5's count.
$ws.Range("F6").Value = 1

# Row 7 (R9 / 1K resistor): give it its own quantity of 1 (previously blank).
$ws.Range("F7").Value = 1

# Move/save the active selection to F7, matching the saved workbook state.
$ws.Range("F7").Select()
